$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.166.69'
$ws.Range('E2').Value = '  +0.71%  '

$ws.Range('D3').Value = '1.901.45'
$ws.Range('E3').Value = '  +1.13%  '

$style = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = $style
$ws.Range('E4').Value = '  +0.58%  '

$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.11'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -0.13%  '

$ws.Range('E6').Value = '  +0.39%  '

$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5256'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  +1.91%  '

$ws.Range('E8').Value = '  +1.67%  '

$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07253'
$ws.Range('D9').Style = $style

$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.16'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  +1.94%  '

$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8987'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  -0.03%  '

$style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08362'
$ws.Range('D12').Style = $style

$ws.Range('D13').Value = '1.897.48'
$ws.Range('E13').Value = '  +0.35%  '

$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.80'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  -0.15%  '

$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.266'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +0.41%  '

$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  +0.59%  '

$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008605'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +1.45%  '

$ws.Range('E18').Value = '  +1.94%  '

$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  +0.37%  '

$ws.Range('D20').Value = '27.209.35'
$ws.Range('E20').Value = '  +0.78%  '

$ws.Range('E21').Value = '  +0.74%  '

$ws.Range('D22').Value = '2.131.63'
$ws.Range('E22').Value = '  +0.78%  '

$ws.Range('E23').Value = '  +1.85%  '

$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.428'
$ws.Range('D24').Style = $style

$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.284'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +7.98%  '

$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '146.47'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  +0.62%  '

$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.759'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  -1.30%  '

$ws.Range('E28').Value = '  +0.68%  '

$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.78'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +0.31%  '

$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.927'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +0.27%  '

$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.783'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  +0.25%  '

$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09247'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  +0.72%  '

$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8115'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  +7.42%  '

$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05052'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  +0.42%  '

$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.239'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  +5.05%  '

$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.958'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -1.31%  '

$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.349'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  +2.16%  '

$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.567'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  +3.15%  '

$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5694'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  +1.97%  '

$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01976'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  -0.77%  '

$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.074'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +0.04%  '

$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.663'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +1.31%  '

$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.940'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  +2.13%  '

$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '118.20'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +1.29%  '

$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1512'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  +0.68%  '

$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4831'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +0.96%  '

$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  +0.35%  '

$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.610'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +2.88%  '

$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.45'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  +1.07%  '

$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.57'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  +0.35%  '
